$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("messages")

# Set distinct markers in A400..A405
for ($i = 400; $i -le 405; $i++) {
  $ws.Range("A$i").Value2 = "ROW_$i"
}

# Shift down by one within 400-406 (simulate insert at 400)
for ($i = 405; $i -ge 400; $i--) {
  $src = $i
  $dst = $i + 1
  $ws.Range("A$src").Copy($ws.Range("A$dst"))
}
$ws.Range("A400").Value2 = "NEW_400"

for ($i = 400; $i -le 406; $i++) {
  $v = $ws.Range("A$i").Value2
  Write-Host "A$i =" 
  Write-Host $v
}
